$d = $word.ActiveDocument

# ===========================================================================
# "Artifact Database Android App" entry
# ===========================================================================

# 1) The date "July 2024 - August 2024" used to be split across three runs
#    ("July 2024 -", " ", "August 2024"); consolidate into a single run
#    (text itself is unchanged).
$rng = $d.Content
$rng.Find.Execute("July 2024 – August 2024", $true, $false, $false, $false, `
                   $false, $true, 1, $false, "", 0) | Out-Null
$rng.Delete()
$rng.InsertAfter("July 2024 – August 2024")

# 2) Add three new "List Bullet" accomplishment lines right after the
#    heading line, before the existing "Implemented a vibrant..." bullet.
$rng = $d.Content
$rng.Find.Execute("July 2024 – August 2024", $true, $false, $false, $false, `
                   $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("Leveraged Firebase and Java to develop a mobile app that efficiently manages and organizes museum artifacts.")
$p = $rng.Paragraphs(1)
$p.Style = "List Bullet"

$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("Led a team of 6 as Scrum Master, facilitating ")
$rng.Collapse(0)
$rng.InsertAfter("standups ")
$rng.Collapse(0)
$rng.InsertAfter("and sprint planning, communicating tasks and deadlines, and ")
$rng.Collapse(0)
$rng.InsertAfter("enabling smooth")
$rng.Collapse(0)
$rng.InsertAfter(" collaboration on GitHub by resolving merge conflicts, earning a peer evaluation score of 5/5.")
$p = $rng.Paragraphs(1)
$p.Style = "List Bullet"

$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("Troubleshooted and resolved bugs for 3 team members, ensuring the artifact display page and report generation feature met requirements and kept the project on schedule.")
$p = $rng.Paragraphs(1)
$p.Style = "List Bullet"

# 3) Add a new bullet right after "Implemented a vibrant, interactive UI...".
$rng2 = $d.Content
$rng2.Find.Execute("Implemented a vibrant, interactive UI by extending standard AndroidX fragments to create custom UI elements.", `
                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.InsertParagraphAfter()
$rng2.Collapse(0)
$rng2.MoveStart(1, 1)
$rng2.InsertAfter("Wrote ")
$rng2.Collapse(0)
$rng2.InsertAfter("maintainable and testable code by utilizing the Model-View-Presenter (similar to MVVM) architecture.")
$p2 = $rng2.Paragraphs(1)
$p2.Style = "List Bullet"

# ===========================================================================
# "MultiClock" entry
# ===========================================================================

# 4) Reword the ending of the "Developed a full stack web app..." bullet.
$rng3 = $d.Content
$rng3.Find.Execute("that allows multiplayer board game players to use a timer.", `
                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Delete()
$rng3.InsertAfter("that serves as a timer for multiplayer board games.")

# 5) Add a new bullet right after that reworded line.
$rng4 = $d.Content
$rng4.Find.Execute("Developed a full stack web app using React for the frontend, Node and Express for the backend, and MongoDB for the database, resulting in an easy-to-use tool that serves as a timer for multiplayer board games.", `
                    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng4.InsertParagraphAfter()
$rng4.Collapse(0)
$rng4.MoveStart(1, 1)
$rng4.InsertAfter("Enabled use of the app on both desktop and mobile by implementing a responsive web design.")
$p4 = $rng4.Paragraphs(1)
$p4.Style = "List Bullet"
